$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 799
$ws.Range("F3").Value = 536
$ws.Range("F4").Value = 283
$ws.Range("F5").Value = 490
$ws.Range("F6").Value = 1132
$ws.Range("F7").Value = 317
$ws.Range("F8").Value = 34
$ws.Range("F9").Value = 113
$ws.Range("F10").Value = 112
$ws.Range("F11").Value = 1148
$ws.Range("F14").Value = 792
$ws.Range("F15").Value = 816
$ws.Range("F16").Value = 185
$ws.Range("F17").Value = 47
$ws.Range("F18").Value = 66
$ws.Range("F20").Value = 190
$ws.Range("F21").Value = 1717
$ws.Range("F22").Value = 2348
$ws.Range("F23").Value = 652
$ws.Range("F24").Value = 68
$ws.Range("F25").Value = 1906
$ws.Range("F26").Value = 334
$ws.Range("F27").Value = 2769
$ws.Range("F28").Value = 510
$ws.Range("F30").Value = 682
$ws.Range("F32").Value = 100
$ws.Range("F33").Value = 95
$ws.Range("F34").Value = 957
$ws.Range("F35").Value = 1691
$ws.Range("F36").Value = 330
$ws.Range("F38").Value = 534
$ws.Range("F39").Value = 155
$ws.Range("F40").Value = 116
$ws.Range("F41").Value = 155
$ws.Range("F42").Value = 14

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G3").Value = 224
$ws.Range("F8").Value = 1
$ws.Range("F12").Value = 71

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 799
$ws.Range("F4").Value = 536
$ws.Range("F5").Value = 283
$ws.Range("F6").Value = 490
$ws.Range("F7").Value = 1132
$ws.Range("F8").Value = 317
$ws.Range("F9").Value = 34
$ws.Range("F10").Value = 113
$ws.Range("F11").Value = 112
$ws.Range("F12").Value = 1148
$ws.Range("F14").Value = 792
$ws.Range("F15").Value = 816
$ws.Range("F16").Value = 185
$ws.Range("G17").Value = 224
$ws.Range("G18").Value = 224
$ws.Range("F20").Value = 47
$ws.Range("F22").Value = 66
$ws.Range("F23").Value = 190
$ws.Range("F24").Value = 1717
$ws.Range("F25").Value = 2348
$ws.Range("F26").Value = 652
$ws.Range("F27").Value = 68
$ws.Range("F30").Value = 2770
$ws.Range("F31").Value = 510
$ws.Range("F32").Value = 1
$ws.Range("F37").Value = 71
$ws.Range("F38").Value = 682
$ws.Range("F40").Value = 100
$ws.Range("F41").Value = 95
$ws.Range("F42").Value = 957
$ws.Range("F43").Value = 1691
$ws.Range("F45").Value = 330
$ws.Range("F46").Value = 534
$ws.Range("F47").Value = 155
$ws.Range("F48").Value = 116
$ws.Range("F49").Value = 155
